$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to remain text formatted so numeric-looking strings are not
# auto-converted to numbers by Excel (matches original inlineStr storage).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "60.429.89"
$ws.Range("E2").Value = "  -2.41%  "

$ws.Range("D3").Value = "3.311.57"
$ws.Range("E3").Value = "  -3.02%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "559.04"
$ws.Range("E5").Value = "  -2.82%  "

$ws.Range("D6").Value = "142.57"
$ws.Range("E6").Value = "  -4.02%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").Value = "3.312.92"
$ws.Range("E8").Value = "  -2.99%  "

$ws.Range("E9").Value = "  -2.67%  "

$ws.Range("D10").Value = "7.87"
$ws.Range("E10").Value = "  -2.34%  "

$ws.Range("E11").Value = "  -3.47%  "

$ws.Range("D12").Value = "0.408"
$ws.Range("E12").Value = "  -1.00%  "

$ws.Range("D13").Value = "3.880.88"
$ws.Range("E13").Value = "  -3.00%  "

$ws.Range("E14").Value = "  +0.29%  "

$ws.Range("D15").Value = "26.90"
$ws.Range("E15").Value = "  -4.75%  "

$ws.Range("D16").Value = "3.311.51"
$ws.Range("E16").Value = "  -3.48%  "

$ws.Range("E17").Value = "  -2.58%  "

$ws.Range("D18").Value = "60.437.83"
$ws.Range("E18").Value = "  -2.43%  "

$ws.Range("D19").Value = "6.14"
$ws.Range("E19").Value = "  -3.03%  "

$ws.Range("D20").Value = "14.33"
$ws.Range("E20").Value = "  -0.57%  "

$ws.Range("D21").Value = "8.68"
$ws.Range("E21").Value = "  -2.00%  "

$ws.Range("D22").Value = "374.75"
$ws.Range("E22").Value = "  -1.54%  "

$ws.Range("D23").Value = "74.87"
$ws.Range("E23").Value = "  -0.34%  "

$ws.Range("E24").Value = "  -0.06%  "

$ws.Range("E25").Value = "  -4.55%  "

$ws.Range("D26").Value = "3.447.84"
$ws.Range("E26").Value = "  -3.05%  "

$ws.Range("E27").Value = "  -7.46%  "

$ws.Range("E28").Value = "  -4.14%  "

$ws.Range("D29").Value = "0.991"
$ws.Range("E29").Value = "  -1.10%  "

$ws.Range("D30").Value = "7.18"
$ws.Range("E30").Value = "  -5.73%  "

$ws.Range("E31").Value = "  -0.04%  "

$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "2.05"
$ws.Range("E32").Value = "  -3.04%  "

$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Value = "7.59"
$ws.Range("E33").Value = "  -3.70%  "

$ws.Range("D34").Value = "22.68"
$ws.Range("E34").Value = "  -1.59%  "

$ws.Range("D35").Value = "1.24"
$ws.Range("E35").Value = "  -7.39%  "

$ws.Range("E36").Value = "  -5.29%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "1.53"
$ws.Range("E37").Value = "  -3.41%  "

$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").Value = "166.35"
$ws.Range("E38").Value = "  -1.91%  "

$ws.Range("E39").Value = "  -2.11%  "

$ws.Range("B40").Value = "RenzoRestakedETH"
$ws.Range("C40").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D40").Value = "3.343.43"
$ws.Range("E40").Value = "  -3.08%  "

$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").Value = "26.85"
$ws.Range("E41").Value = "  -13.04%  "

$ws.Range("D42").Value = "0.0732"
$ws.Range("E42").Value = "  -5.31%  "

$ws.Range("E43").Value = "  -1.33%  "

$ws.Range("D44").Value = "0.753"
$ws.Range("E44").Value = "  -2.87%  "

$ws.Range("D45").Value = "4.17"
$ws.Range("E45").Value = "  -4.46%  "

$ws.Range("D46").Value = "1.58"
$ws.Range("E46").Value = "  -4.95%  "

$ws.Range("E47").Value = "  -3.82%  "

$ws.Range("D48").Value = "2.372.81"
$ws.Range("E48").Value = "  -6.56%  "

$ws.Range("E49").Value = "  -0.05%  "

$ws.Range("D50").Value = "6.43"
$ws.Range("E50").Value = "  -6.54%  "

$ws.Range("D51").Value = "21.28"
$ws.Range("E51").Value = "  -5.00%  "
